# "Primeira aula de AndroidStudio" -- add a new attendance column (C) for
# "Aula 14" / "Grupo 5", marking who showed up, on the single worksheet
# "Planilha1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in column C, mirroring the existing B1/B2 headers
# ("Aula 10" / "Grupo 1") which are centered (style index 1).
$ws.Range("C1").Value = "Aula 14"
$ws.Range("C2").Value = "Grupo 5"

# Attendance marks for this new class date.
# Row 10  -> Juliana de Carvalho Fernandes
# Row 22  -> Monick Hellen Nogueira Macena
# Row 30  -> Roberto Freixeira da Silva Junior (already marked "X" in B30)
$ws.Range("C10").Value = "x"
$ws.Range("C10").HorizontalAlignment = -4108   # xlCenter, matches B-column style

$ws.Range("C22").Value = "x"
$ws.Range("C22").HorizontalAlignment = -4108   # xlCenter

$ws.Range("C30").Value = "X"
$ws.Range("C30").HorizontalAlignment = -4108   # xlCenter

# Restore the view roughly where the author left it (best effort; some
# window-chrome/scroll-position properties are not all settable here).
$ws.Range("E8").Select()
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A4") } catch {}
